$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 249.75  # was 275.5
$ws.Range("I4").Value = 267.6  # was 275.5
$ws.Range("J4").Value = 220  # was 0
$ws.Range("K4").Value = 267.6  # was 275.5
$ws.Range("L4").Value = 220  # was 0
$ws.Range("M4").Value = -153.6  # was -161.5
$ws.Range("N4").Value = -448  # was (empty)
$ws.Range("H10").Value = 5850  # was 10000
$ws.Range("I10").Value = 5050  # was 10000
$ws.Range("J10").Value = 6250  # was 10000
$ws.Range("K10").Value = 5050  # was 10000
$ws.Range("L10").Value = 6250  # was 10000
$ws.Range("M10").Value = -4757  # was -9707
$ws.Range("N10").Value = -6836  # was -10586
$ws.Range("H113").Value = 3981.25  # was 4582.1055
$ws.Range("I113").Value = 2287.5  # was 2708.5715
$ws.Range("K113").Value = 2287.5  # was 2708.5715
$ws.Range("M113").Value = 966.5  # was 545.4285

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 17758  # was 18058.2
$ws.Range("J28").Value = 34172.332  # was 35173
$ws.Range("L28").Value = 34172.332  # was 35173
$ws.Range("N28").Value = -34556.332  # was -35557
$ws.Range("H45").Value = 1766.6666  # was 1228.2
$ws.Range("I45").Value = 1600  # was 1012.3
$ws.Range("J45").Value = 1850  # was 1660
$ws.Range("K45").Value = 1600  # was 1012.3
$ws.Range("L45").Value = 1850  # was 1660
$ws.Range("M45").Value = -1223  # was -635.3
$ws.Range("N45").Value = -2604  # was -2414
$ws.Range("H61").Value = 3760.5334  # was 2749.8684
$ws.Range("I61").Value = 3150  # was 1541.8572
$ws.Range("J61").Value = 4167.5557  # was 4242.1177
$ws.Range("K61").Value = 3150  # was 1541.8572
$ws.Range("L61").Value = 4167.5557  # was 4242.1177
$ws.Range("M61").Value = -2938  # was -1329.8572
$ws.Range("N61").Value = -4591.5557  # was -4666.1177
$ws.Range("H99").Value = 17758  # was 18058.2
$ws.Range("J99").Value = 34172.332  # was 35173
$ws.Range("L99").Value = 34172.332  # was 35173
$ws.Range("N99").Value = -40162.332  # was -41163
$ws.Range("H110").Value = 1706.0869  # was 2027.5333
$ws.Range("I110").Value = 1694  # was 2225
$ws.Range("J110").Value = 1728.75  # was 1801.8572
$ws.Range("K110").Value = 1694  # was 2225
$ws.Range("L110").Value = 1728.75  # was 1801.8572
$ws.Range("M110").Value = 351  # was -180
$ws.Range("N110").Value = -5818.75  # was -5891.8572
$ws.Range("H136").Value = 3760.5334  # was 2749.8684
$ws.Range("I136").Value = 3150  # was 1541.8572
$ws.Range("J136").Value = 4167.5557  # was 4242.1177
$ws.Range("K136").Value = 9450  # was 4625.571599999999
$ws.Range("L136").Value = 12502.6671  # was 12726.3531
$ws.Range("M136").Value = -6900  # was -2075.571599999999
$ws.Range("N136").Value = -17602.6671  # was -17826.3531

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 500  # was 980
$ws.Range("I22").Value = 500  # was 0
$ws.Range("J22").Value = 0  # was 980
$ws.Range("K22").Value = 500  # was 0
$ws.Range("L22").ClearContents()  # was 980
$ws.Range("M22").Value = -327  # was (empty)
$ws.Range("N22").Value = 0  # was -1326
$ws.Range("H87").Value = 19950  # was 20000
$ws.Range("J87").Value = 19950  # was 20000
$ws.Range("L87").Value = 19950  # was 20000
$ws.Range("N87").Value = -22446  # was -22496
$ws.Range("H90").Value = 19950  # was 20000
$ws.Range("J90").Value = 19950  # was 20000
$ws.Range("L90").Value = 59850  # was 60000
$ws.Range("N90").Value = -72330  # was -72480
$ws.Range("H122").Value = 40000  # was 0
$ws.Range("J122").Value = 40000  # was 0
$ws.Range("L122").Value = 40000  # was 0
$ws.Range("N122").Value = -49800  # was (empty)

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 910.1818  # was 944.1111
$ws.Range("I16").Value = 871.2857  # was 919.8
$ws.Range("J16").Value = 978.25  # was 974.5
$ws.Range("K16").Value = 871.2857  # was 919.8
$ws.Range("L16").Value = 978.25  # was 974.5
$ws.Range("M16").Value = -584.2857  # was -632.8
$ws.Range("N16").Value = -1552.25  # was -1548.5
$ws.Range("H22").Value = 380  # was 384.2857
$ws.Range("I22").Value = 126.666664  # was 136.66667
$ws.Range("K22").Value = 126.666664  # was 136.66667
$ws.Range("M22").Value = 223.333336  # was 213.33333
$ws.Range("H31").Value = 2378.5  # was 2983.303
$ws.Range("I31").Value = 1096.3226  # was 1343.3158
$ws.Range("J31").Value = 5028.3335  # was 5209
$ws.Range("K31").Value = 1096.3226  # was 1343.3158
$ws.Range("L31").Value = 5028.3335  # was 5209
$ws.Range("M31").Value = -801.3226  # was -1048.3158
$ws.Range("N31").Value = -5618.3335  # was -5799
$ws.Range("H34").Value = 2378.5  # was 2983.303
$ws.Range("I34").Value = 1096.3226  # was 1343.3158
$ws.Range("J34").Value = 5028.3335  # was 5209
$ws.Range("K34").Value = 1096.3226  # was 1343.3158
$ws.Range("L34").Value = 5028.3335  # was 5209
$ws.Range("M34").Value = -894.3226  # was -1141.3158
$ws.Range("N34").Value = -5432.3335  # was -5613
$ws.Range("H113").Value = 910.1818  # was 944.1111
$ws.Range("I113").Value = 871.2857  # was 919.8
$ws.Range("J113").Value = 978.25  # was 974.5
$ws.Range("K113").Value = 871.2857  # was 919.8
$ws.Range("L113").Value = 978.25  # was 974.5
$ws.Range("M113").Value = 1298.7143  # was 1250.2
$ws.Range("N113").Value = -5318.25  # was -5314.5
$ws.Range("H134").Value = 1650.579  # was 1822.2693
$ws.Range("I134").Value = 1008.3913  # was 1070.8334
$ws.Range("J134").Value = 2635.2666  # was 2466.3572
$ws.Range("K134").Value = 3025.1739  # was 3212.5002
$ws.Range("L134").Value = 7905.7998  # was 7399.071599999999
$ws.Range("M134").Value = -490.1738999999998  # was -677.5001999999999
$ws.Range("N134").Value = -12975.7998  # was -12469.0716

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1082.1025  # was 1089.7949
$ws.Range("J5").Value = 1834.1666  # was 1850.8334
$ws.Range("L5").Value = 5502.4998  # was 5552.5002
$ws.Range("N5").Value = -5726.4998  # was -5776.5002
$ws.Range("H14").Value = 104.611115  # was 1169.75
$ws.Range("I14").Value = 104.611115  # was 1169.75
$ws.Range("K14").Value = 313.833345  # was 3509.25
$ws.Range("M14").Value = -140.833345  # was -3336.25
$ws.Range("H68").Value = 534  # was 2003
$ws.Range("I68").Value = 502  # was 0
$ws.Range("J68").Value = 550  # was 2003
$ws.Range("K68").Value = 1506  # was 0
$ws.Range("L68").Value = 1650  # was 6009
$ws.Range("M68").Value = -695  # was (empty)
$ws.Range("N68").Value = -3272  # was -7631
$ws.Range("H71").Value = 534  # was 2003
$ws.Range("I71").Value = 502  # was 0
$ws.Range("J71").Value = 550  # was 2003
$ws.Range("K71").Value = 4518  # was 0
$ws.Range("L71").Value = 4950  # was 18027
$ws.Range("M71").Value = -462  # was (empty)
$ws.Range("N71").Value = -13062  # was -26139
$ws.Range("H135").Value = 1082.1025  # was 1089.7949
$ws.Range("J135").Value = 1834.1666  # was 1850.8334
$ws.Range("L135").Value = 16507.4994  # was 16657.5006
$ws.Range("N135").Value = -21577.4994  # was -21727.5006

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4959.467  # was 6145.6665
$ws.Range("I102").Value = 2498.8572  # was 3762
$ws.Range("J102").Value = 7112.5  # was 7337.5
$ws.Range("K102").Value = 2498.8572  # was 3762
$ws.Range("L102").Value = 7112.5  # was 7337.5
$ws.Range("M102").Value = -876.8571999999999  # was -2140
$ws.Range("N102").Value = -10356.5  # was -10581.5
$ws.Range("H126").Value = 2628.8823  # was 2858.1333
$ws.Range("I126").Value = 2360.0833  # was 2735.25
$ws.Range("J126").Value = 3274  # was 2998.5715
$ws.Range("K126").Value = 7080.249899999999  # was 8205.75
$ws.Range("L126").Value = 9822  # was 8995.7145
$ws.Range("M126").Value = -4610.249899999999  # was -5735.75
$ws.Range("N126").Value = -14762  # was -13935.7145

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 775.125  # was 758.8823
$ws.Range("I22").Value = 566.6667  # was 566.5
$ws.Range("J22").Value = 900.2  # was 863.8182
$ws.Range("K22").Value = 566.6667  # was 566.5
$ws.Range("L22").Value = 900.2  # was 863.8182
$ws.Range("M22").Value = -271.6667  # was -271.5
$ws.Range("N22").Value = -1490.2  # was -1453.8182
$ws.Range("H27").Value = 775.125  # was 758.8823
$ws.Range("I27").Value = 566.6667  # was 566.5
$ws.Range("J27").Value = 900.2  # was 863.8182
$ws.Range("K27").Value = 566.6667  # was 566.5
$ws.Range("L27").Value = 900.2  # was 863.8182
$ws.Range("M27").Value = -459.6667  # was -459.5
$ws.Range("N27").Value = -1114.2  # was -1077.8182
